# Clean 2015 and 2014 tower data.
# Column A on Sheet1 held a mix of 2014/2016 Excel date-serials where the
# row should have been a clean, consecutive run of 2015 dates. Replace the
# stored serials with the corrected 2015 values (and the trailing row 53
# date) to match the cleaned-up tower count log.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$correctedDates = @{
    2 = 42164
    3 = 42165
    4 = 42166
    5 = 42167
    6 = 42168
    7 = 42169
    8 = 42170
    9 = 42171
    10 = 42172
    11 = 42173
    12 = 42174
    13 = 42175
    14 = 42176
    15 = 42177
    16 = 42178
    17 = 42179
    18 = 42180
    19 = 42181
    20 = 42182
    21 = 42183
    22 = 42184
    23 = 42185
    24 = 42186
    25 = 42187
    26 = 42188
    27 = 42189
    28 = 42190
    29 = 42191
    30 = 42192
    31 = 42193
    32 = 42194
    33 = 42195
    34 = 42196
    35 = 42197
    36 = 42198
    37 = 42199
    38 = 42200
    39 = 42201
    40 = 42202
    41 = 42203
    42 = 42204
    43 = 42205
    44 = 42206
    45 = 42207
    46 = 42208
    47 = 42209
    48 = 42210
    49 = 42211
    50 = 42212
    51 = 42213
    52 = 42214
    53 = 42227
}

foreach ($row in $correctedDates.Keys) {
    $ws.Cells.Item($row, 1).Value = $correctedDates[$row]
}

# Restore the sheet/selection the author left active after the cleanup.
$ws.Activate()
$ws.Range("A46").Select()
